{"js": "// Apply two small text corrections to the report body:\n// 1) \"... alguns est\u00e3o professores em revis\u00e3o ...\" ->\n//    \"... alguns est\u00e3o com professores em revis\u00e3o ...\"\n// 2) \". Para conseguir dar um andamento satisfat\u00f3rio, ...\" ->\n//    \". Para conseguir andamento satisfat\u00f3rio, ...\"\n\nconst body = context.document.body;\n\n// --- Edit 1: insert \"com \" before \"professores em revis\u00e3o\" ---\nconst hits1 = body.search(\"est\u00e3o professores em revis\u00e3o\", { matchCase: true });\nhits1.load(\"text\");\nawait context.sync();\n\nif (hits1.items.length === 0) {\n  throw new Error(\"Target phrase for edit 1 not found\");\n}\n\nhits1.items[0].insertText(\"est\u00e3o com professores em revis\u00e3o\", \"Replace\");\nawait context.sync();\n\n// --- Edit 2: remove \"dar um \" from \"conseguir dar um andamento\" ---\nconst hits2 = body.search(\"conseguir dar um andamento\", { matchCase: true });\nhits2.load(\"text\");\nawait context.sync();\n\nif (hits2.items.length === 0) {\n  throw new Error(\"Target phrase for edit 2 not found\");\n}\n\nhits2.items[0].insertText(\"conseguir andamento\", \"Replace\");\nawait context.sync();\n", "ps1": "# Apply two small text corrections to the report body:\n# 1) \"... alguns est\u00e3o professores em revis\u00e3o ...\" ->\n#    \"... alguns est\u00e3o com professores em revis\u00e3o ...\"\n# 2) \". Para conseguir dar um andamento satisfat\u00f3rio, ...\" ->\n#    \". Para conseguir andamento satisfat\u00f3rio, ...\"\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: insert \"com \" before \"professores em revis\u00e3o\" ---\n$find1 = \"est\u00e3o professores em revis\u00e3o\"\n$replace1 = \"est\u00e3o com professores em revis\u00e3o\"\n$rng1 = $d.Content\n$rng1.Find.ClearFormatting()\n$rng1.Find.Replacement.ClearFormatting()\n$rng1.Find.Text = $find1\n$rng1.Find.Replacement.Text = $replace1\n$rng1.Find.Execute(\n    $find1,    # FindText\n    $false,    # MatchCase\n    $false,    # MatchWholeWord\n    $false,    # MatchWildcards\n    $false,    # MatchSoundsLike\n    $false,    # MatchAllWordForms\n    $true,     # Forward\n    1,         # Wrap (wdFindContinue)\n    $false,    # Format\n    $replace1, # ReplaceWith\n    2          # Replace (wdReplaceAll)\n)\n\n# --- Edit 2: remove \"dar um \" from \"conseguir dar um andamento\" ---\n$find2 = \"conseguir dar um andamento\"\n$replace2 = \"conseguir andamento\"\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Replacement.ClearFormatting()\n$rng2.Find.Text = $find2\n$rng2.Find.Replacement.Text = $replace2\n$rng2.Find.Execute(\n    $find2,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    $replace2,\n    2\n)\n\n$word.ActiveDocument.Saved = $false\n"}
